$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.801.80"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.191.23"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.00"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.23"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "3.189.31"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.96"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "3.713.00"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "3.185.32"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "63.797.61"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.69"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.18"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.84"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.63"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "0.0₃0740"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.74"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0395"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.19"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "399.32"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D44").Value = "2.804.59"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.28"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.85"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.65"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  -0.69%  "
